$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the stray single-space indentation on the "CaseFile caseFile = ..." line
# inside the evalSpring() function body stored in D10 (was a leading single
# space; should match the 4-space indent used by the surrounding lines).
$funcText = "function Boolean evalSpring(String expression, NextPossibleQueuesModel model)`n" +
            "{`n" +
            "    ExpressionParser ep = new SpelExpressionParser();`n" +
            "    Expression exp = ep.parseExpression(expression);`n" +
            "    EvaluationContext ec = new StandardEvaluationContext();`n" +
            " `n" +
            "    CaseFile caseFile = (CaseFile) model.getBusinessObject();`n" +
            "`n" +
            "    Boolean evaluated = exp.getValue(ec, caseFile, Boolean.class);`n" +
            " `n" +
            "    return evaluated;`n" +
            "}"
$ws.Range("D10").Value = $funcText

# Update the CONDITION formula in row 18 from matching the wildcard queue
# name to matching a (non-existent) "no-such-queue" name.
$ws.Range("C18").Value = 'queue.name.equals("no-such-queue")'

# Reflect the author's on-save cursor/viewport position: scrolled back to
# column A (row anchor unchanged at row 5) with C18 as the active selection.
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C18").Select()
